$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- New "Time of day" column (D) ------------------------------------------------
# Row 10 & 11 carry seconds, so Excel keeps them as h:mm:ss; rows 12 & 13 land on an
# exact minute, so Excel's General-to-Time inference drops the seconds (h:mm).
$ws1.Range("D10").Value = 0.52425925925925931
$ws1.Range("D10").NumberFormat = "h:mm:ss"

$ws1.Range("D11").Value = 0.058159722222222217
$ws1.Range("D11").NumberFormat = "h:mm:ss"

$ws1.Range("D12").Value = 0.11527777777777777
$ws1.Range("D12").NumberFormat = "h:mm"

$ws1.Range("D13").Value = 0.375
$ws1.Range("D13").NumberFormat = "h:mm"

# --- New "Date_Time" column (E) = Date (C) + Time_Of_Day (D) ---------------------
$ws1.Range("E10").Formula = "=C10+D10"
$ws1.Range("E11:E13").Formula = "=C11+D11"
$ws1.Range("E10:E13").NumberFormat = "d\-mmm\ hh:mm:ss"

$ws1.Columns.Item(5).AutoFit()

# --- Page setup for the now-printable sheet ---------------------------------------
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Selection / active tab moves to Sheet1 (was "Random") ------------------------
$ws1.Select() | Out-Null
$ws1.Range("E16").Select() | Out-Null
